$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values (E4:H4) - formatting/style unchanged
$ws.Range("E4").Value = 15994
$ws.Range("F4").Value = 13405
$ws.Range("G4").Value = 9689
$ws.Range("H4").Value = 9711

# Update row 5 values (E5:H5) and switch their style to match F4's style
# (drop the bottom border that the old style carried)
$ws.Range("E5").Value = 3676
$ws.Range("F5").Value = 4562
$ws.Range("G5").Value = 4294
$ws.Range("H5").Value = 4707

$ws.Range("E5:H5").Borders.LineStyle = -4142

# Move the current selection to A3
$ws.Range("A3").Select()
